# Apply "Changed HF+ to UOG+ for consistency" edit to S.Table5 worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title cell (A1): HF -> UOG, HF- -> UOG-
$ws.Range("A1").Value = "Supplemental Table 5: Well pad and active wells variation explained among UOG+ samples only for each dataset. Paired datasets were not subject to this analysis, as their respective UOG- samples were excluded."

# Column F (Dataset/measurement labels): HF+ -> UOG+
$ws.Range("F4").Value = "Metatranscriptomics antimicrobial resistance genes UOG+ samples"
$ws.Range("F5").Value = "Metatranscriptomics UOG+ samples"
$ws.Range("F6").Value = "Metatranscriptomics microbial composition (metatranscriptome) UOG+ samples"
$ws.Range("F7").Value = "Sediment 16S rRNA gene UOG+ samples"
$ws.Range("F8").Value = "Water 16S rRNA gene UOG+ samples"

$ws.Range("F9").Value = "Metatranscriptomics antimicrobial resistance genes UOG+ samples"
$ws.Range("F10").Value = "Metatranscriptomics UOG+ samples"
$ws.Range("F11").Value = "Metatranscriptomics microbial composition (metatranscriptome) UOG+ samples"
$ws.Range("F12").Value = "Sediment 16S rRNA gene UOG+ samples"
$ws.Range("F13").Value = "Water 16S rRNA gene UOG+ samples"

# Widen column A to fit content (target stored width ~23.7109375 characters)
$ws.Columns.Item(1).ColumnWidth = 22.8

# Update the active selection to F6, matching the saved view state
$ws.Range("F6").Select()
